$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "1. The system must be able to handle a high volume of simultaneous loan applications without performance degradation."
$ws.Range("A4").Value = "2. The system must ensure the privacy and security of customer data, including Aadhaar/PAN details, credit scores, and KYC documents. Compliance with the Indian IT Act 2000 and the Aadhaar Act 2016 is mandatory."
$ws.Range("A5").Value = "3. The system must have an uptime of at least 99.5% during peak business hours."
$ws.Range("A6").Value = "4. The system must support multi-lingual interface, specifically English and Hindi, as per the Official Languages Act, 1963."
$ws.Range("A7").Value = "5. The system must be able to complete the loan application process within a maximum of 2 hours, reducing the current TAT from 5 days."
$ws.Range("A8").Value = "6. The system must ensure 90% accuracy in eligibility checks, to minimize the risk of loan default."
$ws.Range("A9").Value = "7. The system must ensure real-time verification of KYC documents with an accuracy of 80% to reduce manual data entry."
$ws.Range("A10").Value = "8. The system must be able to integrate with credit bureaus (CIBIL and Experian) for fetching credit scores, complying with the Credit Information Companies (Regulation) Act, 2005."
$ws.Range("A11").Value = "9. The system must ensure secure OTP-based authentication for e-agreement signing and loan disbursement, as per the guidelines of the Reserve Bank of India."
$ws.Range("A12").Value = "10. The system must be able to integrate with core banking for loan disbursement, ensuring secure and seamless transactions."
$ws.Range("A13").Value = "11. The system must be able to support the use of DigiLocker for e-sign integration, as per the guidelines of the Ministry of Electronics and Information Technology (MeitY)."
$ws.Range("A14").Value = "12. The system must ensure that loan disbursement is only through verified bank accounts, as per the Prevention of Money Laundering Act, 2002."
$ws.Range("A15").Value = "13. The system must be ready for UAT by 1st June 2025 and go-live by 1st July 2025, as per the project timeline."
$ws.Range("A16").Value = "14. The system must be adaptable to compliance changes due to upcoming RBI guidelines."
$ws.Range("A17").Value = "15. The system must ensure that the AI-based eligibility and risk assessment module complies with the AI Ethics Guidelines set by the Ministry of Electronics and Information Technology (MeitY)."
